# Applies the "Updated cryptos list" data refresh to Sheet1.
# For each changed cell we set the literal text value coming from the
# source feed. Cells in column D whose new text parses as a plain number
# (e.g. "1.007") are entered with a leading apostrophe so Excel keeps them
# as text (matching the original inlineStr storage) instead of silently
# converting them to numeric values and losing trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.225.59"
$ws.Range("E2").Value = "  -4.12%  "
# Row 3
$ws.Range("D3").Value = "1.655.54"
$ws.Range("E3").Value = "  -3.56%  "
# Row 4
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.21%  "
# Row 5
$ws.Range("D5").Value = "'216.41"
$ws.Range("E5").Value = "  -3.65%  "
# Row 6
$ws.Range("D6").Value = "'0.5134"
$ws.Range("E6").Value = "  -3.06%  "
# Row 7
$ws.Range("D7").Value = "'1.007"
$ws.Range("E7").Value = "  +0.15%  "
# Row 8
$ws.Range("D8").Value = "'0.2601"
$ws.Range("E8").Value = "  -1.97%  "
# Row 9
$ws.Range("D9").Value = "'0.06462"
$ws.Range("E9").Value = "  -3.61%  "
# Row 10
$ws.Range("D10").Value = "'19.94"
$ws.Range("E10").Value = "  -4.62%  "
# Row 11
$ws.Range("D11").Value = "'0.07841"
$ws.Range("E11").Value = "  +2.13%  "
# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.306"
$ws.Range("E12").Value = "  -4.10%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.653.26"
$ws.Range("E13").Value = "  -3.79%  "
# Row 14
$ws.Range("D14").Value = "1.885.61"
$ws.Range("E14").Value = "  -3.42%  "
# Row 15
$ws.Range("D15").Value = "'0.5541"
$ws.Range("E15").Value = "  -4.60%  "
# Row 16
$ws.Range("D16").Value = "0.0₅8047"
$ws.Range("E16").Value = "  -2.00%  "
# Row 17
$ws.Range("D17").Value = "'64.23"
$ws.Range("E17").Value = "  -5.32%  "
# Row 18
$ws.Range("D18").Value = "26.238.48"
$ws.Range("E18").Value = "  -4.16%  "
# Row 19
$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  +0.06%  "
# Row 20
$ws.Range("D20").Value = "'210.93"
$ws.Range("E20").Value = "  -5.13%  "
# Row 21
$ws.Range("D21").Value = "'4.419"
$ws.Range("E21").Value = "  -5.23%  "
# Row 22
$ws.Range("D22").Value = "'10.09"
$ws.Range("E22").Value = "  -3.50%  "
# Row 23
$ws.Range("D23").Value = "'6.030"
$ws.Range("E23").Value = "  +0.24%  "
# Row 24
$ws.Range("E24").Value = "  +0.16%  "
# Row 25
$ws.Range("D25").Value = "'144.83"
$ws.Range("E25").Value = "  -0.34%  "
# Row 26
$ws.Range("D26").Value = "'1.792"
$ws.Range("E26").Value = "  +4.96%  "
# Row 27
$ws.Range("D27").Value = "'0.1179"
$ws.Range("E27").Value = "  -2.37%  "
# Row 28
$ws.Range("D28").Value = "'7.023"
$ws.Range("E28").Value = "  -3.13%  "
# Row 29
$ws.Range("E29").Value = "  -2.16%  "
# Row 30
$ws.Range("D30").Value = "'0.05105"
$ws.Range("E30").Value = "  -5.31%  "
# Row 31
$ws.Range("D31").Value = "'1.243"
$ws.Range("E31").Value = "  -4.01%  "
# Row 32
$ws.Range("D32").Value = "'3.367"
$ws.Range("E32").Value = "  -3.31%  "
# Row 33
$ws.Range("D33").Value = "'3.233"
$ws.Range("E33").Value = "  -5.41%  "
# Row 34
$ws.Range("D34").Value = "'1.566"
$ws.Range("E34").Value = "  -4.42%  "
# Row 35
$ws.Range("D35").Value = "'2.732"
$ws.Range("E35").Value = "  -4.59%  "
# Row 36
$ws.Range("D36").Value = "'0.9259"
$ws.Range("E36").Value = "  -2.87%  "
# Row 37
$ws.Range("D37").Value = "'2.355"
$ws.Range("E37").Value = "  -1.74%  "
# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.5740"
$ws.Range("E38").Value = "  -2.66%  "
# Row 39
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.166.78"
$ws.Range("E39").Value = "  +1.44%  "
# Row 40
$ws.Range("D40").Value = "'0.01591"
$ws.Range("E40").Value = "  -3.73%  "
# Row 41
$ws.Range("D41").Value = "'2.558"
$ws.Range("E41").Value = "  -0.18%  "
# Row 42
$ws.Range("D42").Value = "'1.007"
$ws.Range("E42").Value = "  +0.08%  "
# Row 43
$ws.Range("D43").Value = "'5.719"
$ws.Range("E43").Value = "  -2.15%  "
# Row 44
$ws.Range("D44").Value = "'0.8251"
$ws.Range("E44").Value = "  -1.98%  "
# Row 45
$ws.Range("D45").Value = "'100.37"
$ws.Range("E45").Value = "  -0.68%  "
# Row 46
$ws.Range("D46").Value = "1.797.87"
$ws.Range("E46").Value = "  -3.31%  "
# Row 47
$ws.Range("E47").Value = "  -0.31%  "
# Row 48
$ws.Range("D48").Value = "'0.4549"
$ws.Range("E48").Value = "  -0.73%  "
# Row 49
$ws.Range("D49").Value = "'55.47"
$ws.Range("E49").Value = "  -4.15%  "
# Row 50
$ws.Range("D50").Value = "'1.007"
$ws.Range("E50").Value = "  +0.32%  "
# Row 51
$ws.Range("D51").Value = "'7.877"
$ws.Range("E51").Value = "  -3.05%  "
